$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("X1").Value = 0.86680005771130819
$ws.Range("A2").Value = 0.69604619271405221
$ws.Range("H2").Value = 0.85502337726983435
$ws.Range("AF2").Value = 0.63356994946955991
$ws.Range("B3").Value = 0.83358582451646068
$ws.Range("BF3").Value = 0.93015753473127494
$ws.Range("E4").Value = 0.82136007281672019
$ws.Range("F4").Value = 0.76604914566636007
$ws.Range("R4").Value = 0.89695165709612179
$ws.Range("C5").Value = 0.82693899910519808
$ws.Range("G5").Value = 0.87092158808039843
$ws.Range("F7").Value = 0.86996365047366364
$ws.Range("V7").Value = 0.89301660332755051
$ws.Range("AR7").Value = 0.91542665410659085
$ws.Range("F8").Value = 0.67286641743535769
$ws.Range("G8").Value = 0.81221904022357183
$ws.Range("J8").Value = 0.56088634180341712
$ws.Range("K9").Value = 0.90390105047462677
$ws.Range("S9").Value = 0.98729729920209852
$ws.Range("BE9").Value = 0.89111790748697661
$ws.Range("K10").Value = 0.77818739494329492
$ws.Range("L10").Value = 0.6695234798245463
$ws.Range("BI10").Value = 0.99907847164827923
$ws.Range("M11").Value = 0.93839780892932068
$ws.Range("K12").Value = 0.82069711566996606
$ws.Range("M12").Value = 0.97943287626866926
$ws.Range("N12").Value = 0.79259855929036416
$ws.Range("N13").Value = 0.86461527506230107
$ws.Range("N15").Value = 0.69552727634802458
$ws.Range("K16").Value = 0.93712924771919781
$ws.Range("O16").Value = 0.89981392093418022
$ws.Range("O17").Value = 0.9992203561512869
$ws.Range("BP17").Value = 0.99704030991487302
$ws.Range("Q18").Value = 0.86797847515473037
$ws.Range("S18").Value = 0.79902137710216803
$ws.Range("Q19").Value = 0.74567109292961375
$ws.Range("U20").Value = 0.66454790212967851
$ws.Range("V20").Value = 0.85869606979207658
$ws.Range("S21").Value = 0.85968806183419488
$ws.Range("V21").Value = 0.92202535980502986
$ws.Range("BD21").Value = 0.73238883902594765
$ws.Range("U23").Value = 0.74755336441270304
$ws.Range("V23").Value = 0.66564878463380439
$ws.Range("X23").Value = 0.60014501832584366
$ws.Range("V24").Value = 0.69186241970190476
$ws.Range("Z24").Value = 0.84666908564722898
$ws.Range("W25").Value = 0.92349067531920515
$ws.Range("AN25").Value = 0.88187739809079446
$ws.Range("AB26").Value = 0.82322780238152238
$ws.Range("Y27").Value = 0.82741919204961545
$ws.Range("AC27").Value = 0.63463249269747024
$ws.Range("AD28").Value = 0.9469775601681405
$ws.Range("AW28").Value = 0.96266807781868713
$ws.Range("AC30").Value = 0.96547989442965321
$ws.Range("AQ30").Value = 0.7473915256065764
$ws.Range("AC31").Value = 0.91678596825273528
$ws.Range("AD31").Value = 0.92049689854485162
$ws.Range("AG31").Value = 0.9515176482015959
$ws.Range("AH31").Value = 0.81971020886410229
$ws.Range("AG32").Value = 0.60497498378976267
$ws.Range("AH33").Value = 0.90747844434096703
$ws.Range("N34").Value = 0.74585664321817002
$ws.Range("AI34").Value = 0.65169301412717806
$ws.Range("AG35").Value = 0.62436690107173332
$ws.Range("AK35").Value = 0.97823678767093125
$ws.Range("AH36").Value = 0.91422614695856907
$ws.Range("AI36").Value = 0.89557572224625615
$ws.Range("AK36").Value = 0.97825982982958726
$ws.Range("AL36").Value = 0.99909555518300741
$ws.Range("AO36").Value = 0.88337443437741947
$ws.Range("AK38").Value = 0.86942142467968531
$ws.Range("AM38").Value = 0.88264823241066503
$ws.Range("BE38").Value = 0.95452103075477268
$ws.Range("AK39").Value = 0.93852534751202388
$ws.Range("AO39").Value = 0.98188470750870727
$ws.Range("H40").Value = 0.87062236908449298
$ws.Range("I40").Value = 0.87339060703256099
$ws.Range("AM40").Value = 0.78771311754234752
$ws.Range("AP40").Value = 0.62919235937965556
$ws.Range("BL42").Value = 0.97703163809897142
$ws.Range("AO43").Value = 0.83464345196514644
$ws.Range("AP43").Value = 0.99076363255672062
$ws.Range("AS43").Value = 0.68212773165301621
$ws.Range("AQ44").Value = 0.55638049389249811
$ws.Range("AU44").Value = 0.98231708225056424
$ws.Range("Z46").Value = 0.76356726158485455
$ws.Range("AS46").Value = 0.93919883510019797
$ws.Range("N47").Value = 0.72245127747316118
$ws.Range("AS47").Value = 0.78165261169019007
$ws.Range("AT47").Value = 0.8598425174210873
$ws.Range("AW48").Value = 0.97390665865915871
$ws.Range("AU49").Value = 0.84847107006672207
$ws.Range("AV50").Value = 0.92843982473813824
$ws.Range("AW50").Value = 0.88610013234555463
$ws.Range("AZ50").Value = 0.850898948747427
$ws.Range("AX51").Value = 0.74820508110515827
$ws.Range("BA51").Value = 0.77429619665521487
$ws.Range("BG51").Value = 0.68123957050903194
$ws.Range("AU52").Value = 0.70353324135139739
$ws.Range("AY52").Value = 0.6550964807740054
$ws.Range("BA52").Value = 0.93756673963542647
$ws.Range("BC53").Value = 0.93288511223680159
$ws.Range("X54").Value = 0.78663608388853634
$ws.Range("AT54").Value = 0.98093829035708302
$ws.Range("BA54").Value = 0.95514933608278674
$ws.Range("P55").Value = 0.90593112639355466
$ws.Range("BB55").Value = 0.95632729556581419
$ws.Range("BD55").Value = 0.73428628273259733
$ws.Range("BE55").Value = 0.71790650763251085
$ws.Range("AA56").Value = 0.89513439418559693
$ws.Range("Q57").Value = 0.77874491786133426
$ws.Range("BD57").Value = 0.73477018209880041
$ws.Range("E59").Value = 0.82057900570917841
$ws.Range("BE59").Value = 0.72246500378863954
$ws.Range("BF59").Value = 0.99472040525943006
$ws.Range("BH59").Value = 0.85908312992331459
$ws.Range("BF60").Value = 0.94751865014665282
$ws.Range("BL60").Value = 0.84974095107609227
$ws.Range("BK61").Value = 0.92929144096838878
$ws.Range("BH62").Value = 0.77344928983775318
$ws.Range("BL62").Value = 0.78270718928507144
$ws.Range("K63").Value = 0.97543380256655587
$ws.Range("M63").Value = 0.81592651595090016
$ws.Range("AB63").Value = 0.58377265973196624
$ws.Range("BJ63").Value = 0.87780793233086385
$ws.Range("BF64").Value = 0.76009391638915336
$ws.Range("BK64").Value = 0.92961055834890771
$ws.Range("BN65").Value = 0.61074501410846294
$ws.Range("AX66").Value = 0.79254851344871935
$ws.Range("BB66").Value = 0.87853115956875638
$ws.Range("BL66").Value = 0.73858627289362289
$ws.Range("BP66").Value = 0.67685378597132817
$ws.Range("AE67").Value = 0.69345246963604046
$ws.Range("BM67").Value = 0.65446182656322116
$ws.Range("BP67").Value = 0.84529617278590918
$ws.Range("C68").Value = 0.74220130637720083
